# Anonymize the staff roster: rename the "team" column to "comment" and
# replace real staff names with generic "スタッフN" placeholders, moving the
# old team label (A/B) into descriptive free-text shift comments. Also
# drops the red highlight font that used to sit on the C:G shift-symbol
# cells for every staff row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column header: チーム -> コメント
$ws.Range("B2").Value = "コメント"

# Row 4..33: staff name (A) + comment (B)
$rows = @{
    4  = @("スタッフ1", "")
    5  = @("スタッフ2", "夜勤4回まで")
    6  = @("スタッフ3", "夜勤4回まで")
    7  = @("スタッフ4", "")
    8  = @("スタッフ5", "")
    9  = @("スタッフ6", "")
    10 = @("スタッフ7", "")
    11 = @("スタッフ8", "")
    12 = @("スタッフ9", "夜勤3回まで")
    13 = @("スタッフ10", "")
    14 = @("スタッフ11", "")
    15 = @("スタッフ12", "")
    16 = @("スタッフ13", "")
    17 = @("スタッフ14", "")
    18 = @("スタッフ15", "")
    19 = @("スタッフ16", "新人　月前半長夜勤なし　")
    20 = @("スタッフ17", "")
    21 = @("スタッフ18", "")
    22 = @("スタッフ19", "")
    23 = @("スタッフ20", "土日休み日勤のみ")
    24 = @("スタッフ21", "夜勤土日のみ3回まで")
    25 = @("スタッフ22", "長入明　水木金3回まで")
    26 = @("スタッフ23", "")
    27 = @("スタッフ24", "")
    28 = @("スタッフ25", "")
    29 = @("スタッフ26", "")
    30 = @("スタッフ27", "")
    31 = @("スタッフ28", "")
    32 = @("スタッフ29", "")
    33 = @("スタッフ30", "新人　月前半長夜勤なし")
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Cells.Item($r, 1).Value = $vals[0]
    $ws.Cells.Item($r, 2).Value = $vals[1]
}

# Drop the red font on the shift-symbol block (columns C:G) for every staff
# row -- back to plain automatic/black text.
$ws.Range("C4:G33").Font.Color = 0
